$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 166
$ws1.Range("G3").Value = 70
$ws1.Range("F4").Value = 616
$ws1.Range("F5").Value = 3078
$ws1.Range("G5").Value = 70
$ws1.Range("F6").Value = 821
$ws1.Range("G6").Value = 70
$ws1.Range("F9").Value = 467
$ws1.Range("F12").Value = 594
$ws1.Range("F14").Value = 2168
$ws1.Range("F19").Value = 2691
$ws1.Range("F25").Value = 700
$ws1.Range("F26").Value = 700
$ws1.Range("F27").Value = 24
$ws1.Range("F28").Value = 28
$ws1.Range("F30").Value = 30
$ws1.Range("F32").Value = 569
$ws1.Range("F35").Value = 915
$ws1.Range("F36").Value = 4725
$ws1.Range("F37").Value = 288
$ws1.Range("F38").Value = 53
$ws1.Range("F39").Value = 20

# ---- Sheet 2: 演出 ----
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("G2").Value = "不可售"
$ws2.Range("F3").Value = 8
$ws2.Range("F8").Value = 368
$ws2.Range("F25").Value = 6
$ws2.Range("F26").Value = 310
$ws2.Range("F28").Value = 320
$ws2.Range("F32").Value = 36
$ws2.Range("F38").Value = 616

# ---- Sheet 3: 本地生活 ----
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F4").Value = 1476
$ws3.Range("F5").Value = 584
$ws3.Range("F6").Value = 301
$ws3.Range("F7").Value = 287

# ---- Sheet 4: 全部类型 ----
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F3").Value = 1476
$ws4.Range("F4").Value = 584
$ws4.Range("F5").Value = 166
$ws4.Range("F6").Value = 301
$ws4.Range("G7").Value = 70
$ws4.Range("F8").Value = 616
$ws4.Range("F9").Value = 3078
$ws4.Range("G9").Value = 70
$ws4.Range("F10").Value = 821
$ws4.Range("G10").Value = 70
$ws4.Range("F13").Value = 467
$ws4.Range("F15").Value = 8
$ws4.Range("F17").Value = 594
$ws4.Range("F18").Value = 368
$ws4.Range("F21").Value = 2168
$ws4.Range("F27").Value = 2691
$ws4.Range("F33").Value = 287
$ws4.Range("F35").Value = 700
$ws4.Range("F36").Value = 700
$ws4.Range("F37").Value = 24
$ws4.Range("F39").Value = 30
$ws4.Range("F41").Value = 569
$ws4.Range("F42").Value = 310
$ws4.Range("F46").Value = 915
$ws4.Range("F47").Value = 4725
$ws4.Range("F48").Value = 288
$ws4.Range("F49").Value = 53
$ws4.Range("F50").Value = 616
$ws4.Range("F51").Value = 616
